# Edit script: insert new price records for Plátano (Mercado Mayorista Lo Valledor de Santiago)
# and shift existing rows 1143-1228 down by 3 to rows 1146-1231, adding 3 brand-new rows
# at 1143-1145 (date 44746) and 3 more new rows at 1229-1231 (date 44194, carried from the old tail).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(1143, 4).Value = 44746
$ws.Cells.Item(1143, 11).Value = 'Sin especificar'
$ws.Cells.Item(1143, 12).Value = 'Pintón'
$ws.Cells.Item(1143, 13).Value = 1135
$ws.Cells.Item(1143, 15).Value = 21000
$ws.Cells.Item(1143, 16).Value = 20269
$ws.Cells.Item(1143, 19).Value = 1013
$ws.Cells.Item(1144, 4).Value = 44746
$ws.Cells.Item(1144, 12).Value = 'Primera Maduro'
$ws.Cells.Item(1144, 14).Value = 20000
$ws.Cells.Item(1144, 15).Value = 21000
$ws.Cells.Item(1144, 16).Value = 20438
$ws.Cells.Item(1144, 19).Value = 1022
$ws.Cells.Item(1145, 4).Value = 44746
$ws.Cells.Item(1145, 12).Value = 'Primera Pintón'
$ws.Cells.Item(1145, 13).Value = 3175
$ws.Cells.Item(1145, 14).Value = 21000
$ws.Cells.Item(1145, 15).Value = 23000
$ws.Cells.Item(1145, 16).Value = 21706
$ws.Cells.Item(1145, 19).Value = 1085
$ws.Cells.Item(1146, 11).Value = 'Barraganete'
$ws.Cells.Item(1146, 12).Value = 'Primera'
$ws.Cells.Item(1146, 13).Value = 216
$ws.Cells.Item(1146, 14).Value = 19000
$ws.Cells.Item(1146, 15).Value = 20000
$ws.Cells.Item(1146, 16).Value = 19500
$ws.Cells.Item(1146, 19).Value = 975
$ws.Cells.Item(1147, 4).Value = 44386
$ws.Cells.Item(1147, 11).Value = 'Sin especificar'
$ws.Cells.Item(1147, 12).Value = 'Pintón'
$ws.Cells.Item(1147, 13).Value = 640
$ws.Cells.Item(1147, 14).Value = 10000
$ws.Cells.Item(1147, 15).Value = 11000
$ws.Cells.Item(1147, 16).Value = 10500
$ws.Cells.Item(1147, 19).Value = 525
$ws.Cells.Item(1148, 4).Value = 44386
$ws.Cells.Item(1148, 12).Value = 'Primera Maduro'
$ws.Cells.Item(1148, 14).Value = 12000
$ws.Cells.Item(1148, 15).Value = 12000
$ws.Cells.Item(1148, 16).Value = 12000
$ws.Cells.Item(1148, 19).Value = 600
$ws.Cells.Item(1149, 4).Value = 44386
$ws.Cells.Item(1149, 12).Value = 'Primera Pintón'
$ws.Cells.Item(1149, 13).Value = 2130
$ws.Cells.Item(1149, 14).Value = 12000
$ws.Cells.Item(1149, 15).Value = 13000
$ws.Cells.Item(1149, 16).Value = 12549
$ws.Cells.Item(1149, 19).Value = 627
$ws.Cells.Item(1150, 11).Value = 'Barraganete'
$ws.Cells.Item(1150, 12).Value = 'Primera'
$ws.Cells.Item(1150, 13).Value = 324
$ws.Cells.Item(1150, 14).Value = 19000
$ws.Cells.Item(1150, 15).Value = 20000
$ws.Cells.Item(1150, 16).Value = 19500
$ws.Cells.Item(1150, 19).Value = 975
$ws.Cells.Item(1151, 4).Value = 44690
$ws.Cells.Item(1151, 12).Value = 'Maduro'
$ws.Cells.Item(1151, 13).Value = 480
$ws.Cells.Item(1151, 14).Value = 7500
$ws.Cells.Item(1151, 15).Value = 8000
$ws.Cells.Item(1151, 16).Value = 7750
$ws.Cells.Item(1151, 19).Value = 388
$ws.Cells.Item(1152, 4).Value = 44690
$ws.Cells.Item(1152, 12).Value = 'Pintón'
$ws.Cells.Item(1152, 13).Value = 1680
$ws.Cells.Item(1152, 14).Value = 9000
$ws.Cells.Item(1152, 15).Value = 10000
$ws.Cells.Item(1152, 16).Value = 9357
$ws.Cells.Item(1152, 19).Value = 468
$ws.Cells.Item(1153, 4).Value = 44690
$ws.Cells.Item(1153, 12).Value = 'Primera Pintón'
$ws.Cells.Item(1153, 13).Value = 2640
$ws.Cells.Item(1153, 14).Value = 12000
$ws.Cells.Item(1153, 15).Value = 13000
$ws.Cells.Item(1153, 16).Value = 12409
$ws.Cells.Item(1153, 19).Value = 620
$ws.Cells.Item(1154, 4).Value = 44631
$ws.Cells.Item(1154, 13).Value = 1440
$ws.Cells.Item(1154, 14).Value = 15000
$ws.Cells.Item(1154, 15).Value = 16000
$ws.Cells.Item(1154, 16).Value = 15500
$ws.Cells.Item(1154, 19).Value = 775
$ws.Cells.Item(1155, 4).Value = 44631
$ws.Cells.Item(1155, 13).Value = 2480
$ws.Cells.Item(1155, 14).Value = 16000
$ws.Cells.Item(1155, 15).Value = 17000
$ws.Cells.Item(1155, 16).Value = 16500
$ws.Cells.Item(1155, 19).Value = 825
$ws.Cells.Item(1156, 4).Value = 44235
$ws.Cells.Item(1156, 13).Value = 300
$ws.Cells.Item(1156, 15).Value = 10000
$ws.Cells.Item(1156, 16).Value = 10000
$ws.Cells.Item(1156, 19).Value = 500
$ws.Cells.Item(1157, 4).Value = 44235
$ws.Cells.Item(1157, 13).Value = 300
$ws.Cells.Item(1157, 14).Value = 13000
$ws.Cells.Item(1157, 15).Value = 13000
$ws.Cells.Item(1157, 16).Value = 13000
$ws.Cells.Item(1157, 19).Value = 650
$ws.Cells.Item(1158, 4).Value = 44235
$ws.Cells.Item(1158, 12).Value = 'Primera Pintón'
$ws.Cells.Item(1158, 13).Value = 3680
$ws.Cells.Item(1158, 14).Value = 14000
$ws.Cells.Item(1158, 15).Value = 15500
$ws.Cells.Item(1158, 16).Value = 14728
$ws.Cells.Item(1158, 19).Value = 736
$ws.Cells.Item(1159, 12).Value = 'Maduro'
$ws.Cells.Item(1159, 13).Value = 2000
$ws.Cells.Item(1159, 14).Value = 10000
$ws.Cells.Item(1159, 15).Value = 11000
$ws.Cells.Item(1159, 16).Value = 10580
$ws.Cells.Item(1159, 19).Value = 529
$ws.Cells.Item(1160, 4).Value = 44582
$ws.Cells.Item(1160, 11).Value = 'Sin especificar'
$ws.Cells.Item(1160, 12).Value = 'Pintón'
$ws.Cells.Item(1160, 13).Value = 270
$ws.Cells.Item(1160, 14).Value = 12000
$ws.Cells.Item(1160, 15).Value = 12000
$ws.Cells.Item(1160, 16).Value = 12000
$ws.Cells.Item(1160, 19).Value = 600
$ws.Cells.Item(1161, 4).Value = 44582
$ws.Cells.Item(1161, 12).Value = 'Primera Maduro'
$ws.Cells.Item(1161, 13).Value = 1100
$ws.Cells.Item(1161, 14).Value = 13000
$ws.Cells.Item(1161, 15).Value = 14000
$ws.Cells.Item(1161, 16).Value = 13564
$ws.Cells.Item(1161, 19).Value = 678
$ws.Cells.Item(1162, 4).Value = 44582
$ws.Cells.Item(1162, 12).Value = 'Primera Pintón'
$ws.Cells.Item(1162, 13).Value = 3430
$ws.Cells.Item(1162, 14).Value = 13000
$ws.Cells.Item(1162, 15).Value = 14000
$ws.Cells.Item(1162, 16).Value = 13580
$ws.Cells.Item(1162, 19).Value = 679
$ws.Cells.Item(1163, 11).Value = 'Barraganete'
$ws.Cells.Item(1163, 12).Value = 'Primera'
$ws.Cells.Item(1163, 13).Value = 216
$ws.Cells.Item(1163, 14).Value = 15000
$ws.Cells.Item(1163, 15).Value = 16000
$ws.Cells.Item(1163, 16).Value = 15500
$ws.Cells.Item(1163, 19).Value = 775
$ws.Cells.Item(1164, 4).Value = 44307
$ws.Cells.Item(1164, 11).Value = 'Sin especificar'
$ws.Cells.Item(1164, 12).Value = 'Maduro'
$ws.Cells.Item(1164, 13).Value = 300
$ws.Cells.Item(1164, 14).Value = 8000
$ws.Cells.Item(1164, 15).Value = 8000
$ws.Cells.Item(1164, 16).Value = 8000
$ws.Cells.Item(1164, 19).Value = 400
$ws.Cells.Item(1165, 4).Value = 44307
$ws.Cells.Item(1165, 12).Value = 'Pintón'
$ws.Cells.Item(1165, 13).Value = 1500
$ws.Cells.Item(1165, 16).Value = 10400
$ws.Cells.Item(1165, 19).Value = 520
$ws.Cells.Item(1166, 4).Value = 44307
$ws.Cells.Item(1166, 12).Value = 'Primera Pintón'
$ws.Cells.Item(1166, 13).Value = 1980
$ws.Cells.Item(1166, 14).Value = 12000
$ws.Cells.Item(1166, 15).Value = 13000
$ws.Cells.Item(1166, 16).Value = 12424
$ws.Cells.Item(1166, 19).Value = 621
$ws.Cells.Item(1167, 11).Value = 'Barraganete'
$ws.Cells.Item(1167, 12).Value = 'Primera'
$ws.Cells.Item(1167, 13).Value = 216
$ws.Cells.Item(1167, 14).Value = 20000
$ws.Cells.Item(1167, 15).Value = 21000
$ws.Cells.Item(1167, 16).Value = 20500
$ws.Cells.Item(1167, 19).Value = 1025
$ws.Cells.Item(1168, 12).Value = 'Maduro'
$ws.Cells.Item(1168, 13).Value = 880
$ws.Cells.Item(1168, 14).Value = 10000
$ws.Cells.Item(1168, 15).Value = 11000
$ws.Cells.Item(1168, 16).Value = 10545
$ws.Cells.Item(1168, 19).Value = 527
$ws.Cells.Item(1169, 4).Value = 44672
$ws.Cells.Item(1169, 13).Value = 2060
$ws.Cells.Item(1169, 14).Value = 10000
$ws.Cells.Item(1169, 15).Value = 12000
$ws.Cells.Item(1169, 16).Value = 10942
$ws.Cells.Item(1169, 19).Value = 547
$ws.Cells.Item(1170, 4).Value = 44672
$ws.Cells.Item(1170, 12).Value = 'Primera Maduro'
$ws.Cells.Item(1170, 13).Value = 1000
$ws.Cells.Item(1170, 14).Value = 12000
$ws.Cells.Item(1170, 15).Value = 13000
$ws.Cells.Item(1170, 16).Value = 12560
$ws.Cells.Item(1170, 19).Value = 628
$ws.Cells.Item(1171, 4).Value = 44672
$ws.Cells.Item(1171, 11).Value = 'Sin especificar'
$ws.Cells.Item(1171, 12).Value = 'Primera Pintón'
$ws.Cells.Item(1171, 13).Value = 3900
$ws.Cells.Item(1171, 14).Value = 13000
$ws.Cells.Item(1171, 15).Value = 14000
$ws.Cells.Item(1171, 16).Value = 13508
$ws.Cells.Item(1171, 19).Value = 675
$ws.Cells.Item(1172, 4).Value = 44344
$ws.Cells.Item(1172, 11).Value = 'Sin especificar'
$ws.Cells.Item(1172, 12).Value = 'Pintón'
$ws.Cells.Item(1172, 13).Value = 1120
$ws.Cells.Item(1172, 14).Value = 9000
$ws.Cells.Item(1172, 15).Value = 10000
$ws.Cells.Item(1172, 16).Value = 9500
$ws.Cells.Item(1172, 19).Value = 475
$ws.Cells.Item(1173, 4).Value = 44344
$ws.Cells.Item(1173, 11).Value = 'Sin especificar'
$ws.Cells.Item(1173, 13).Value = 1500
$ws.Cells.Item(1173, 14).Value = 10000
$ws.Cells.Item(1173, 15).Value = 12000
$ws.Cells.Item(1173, 16).Value = 11033
$ws.Cells.Item(1173, 19).Value = 552
$ws.Cells.Item(1174, 11).Value = 'Barraganete'
$ws.Cells.Item(1174, 13).Value = 150
$ws.Cells.Item(1174, 14).Value = 12000
$ws.Cells.Item(1174, 15).Value = 12000
$ws.Cells.Item(1174, 16).Value = 12000
$ws.Cells.Item(1174, 19).Value = 600
$ws.Cells.Item(1175, 11).Value = 'Barraganete'
$ws.Cells.Item(1175, 12).Value = 'Primera'
$ws.Cells.Item(1175, 13).Value = 296
$ws.Cells.Item(1175, 14).Value = 20000
$ws.Cells.Item(1175, 15).Value = 21000
$ws.Cells.Item(1175, 16).Value = 20365
$ws.Cells.Item(1175, 19).Value = 1018
$ws.Cells.Item(1176, 11).Value = 'Barraganete'
$ws.Cells.Item(1176, 12).Value = 'Primera Pintón'
$ws.Cells.Item(1176, 13).Value = 150
$ws.Cells.Item(1176, 14).Value = 15000
$ws.Cells.Item(1176, 15).Value = 15000
$ws.Cells.Item(1176, 16).Value = 15000
$ws.Cells.Item(1176, 18).Value = 'Ecuador'
$ws.Cells.Item(1176, 19).Value = 750
$ws.Cells.Item(1177, 12).Value = 'Maduro'
$ws.Cells.Item(1177, 13).Value = 360
$ws.Cells.Item(1177, 15).Value = 13000
$ws.Cells.Item(1177, 16).Value = 13000
$ws.Cells.Item(1177, 18).Value = 'Ecuador'
$ws.Cells.Item(1177, 19).Value = 650
$ws.Cells.Item(1178, 12).Value = 'Pintón'
$ws.Cells.Item(1178, 13).Value = 400
$ws.Cells.Item(1178, 14).Value = 13500
$ws.Cells.Item(1178, 15).Value = 14000
$ws.Cells.Item(1178, 16).Value = 13750
$ws.Cells.Item(1178, 19).Value = 688
$ws.Cells.Item(1179, 4).Value = 44433
$ws.Cells.Item(1179, 12).Value = 'Primera Maduro'
$ws.Cells.Item(1179, 13).Value = 680
$ws.Cells.Item(1179, 14).Value = 13000
$ws.Cells.Item(1179, 15).Value = 13000
$ws.Cells.Item(1179, 16).Value = 13000
$ws.Cells.Item(1179, 18).Value = 'Bolivia'
$ws.Cells.Item(1179, 19).Value = 650
$ws.Cells.Item(1180, 4).Value = 44433
$ws.Cells.Item(1180, 12).Value = 'Primera Pintón'
$ws.Cells.Item(1180, 13).Value = 1330
$ws.Cells.Item(1180, 14).Value = 13000
$ws.Cells.Item(1180, 15).Value = 14000
$ws.Cells.Item(1180, 16).Value = 13421
$ws.Cells.Item(1180, 18).Value = 'Bolivia'
$ws.Cells.Item(1180, 19).Value = 671
$ws.Cells.Item(1181, 4).Value = 44433
$ws.Cells.Item(1181, 12).Value = 'Primera Pintón'
$ws.Cells.Item(1181, 13).Value = 960
$ws.Cells.Item(1181, 14).Value = 15000
$ws.Cells.Item(1181, 15).Value = 16000
$ws.Cells.Item(1181, 16).Value = 15500
$ws.Cells.Item(1181, 19).Value = 775
$ws.Cells.Item(1182, 12).Value = 'Maduro'
$ws.Cells.Item(1182, 13).Value = 1000
$ws.Cells.Item(1182, 14).Value = 7000
$ws.Cells.Item(1182, 15).Value = 7500
$ws.Cells.Item(1182, 16).Value = 7260
$ws.Cells.Item(1182, 19).Value = 363
$ws.Cells.Item(1183, 4).Value = 44707
$ws.Cells.Item(1183, 13).Value = 2220
$ws.Cells.Item(1183, 14).Value = 7000
$ws.Cells.Item(1183, 15).Value = 8000
$ws.Cells.Item(1183, 16).Value = 7523
$ws.Cells.Item(1183, 19).Value = 376
$ws.Cells.Item(1184, 4).Value = 44707
$ws.Cells.Item(1184, 12).Value = 'Primera Maduro'
$ws.Cells.Item(1184, 13).Value = 2000
$ws.Cells.Item(1184, 14).Value = 9000
$ws.Cells.Item(1184, 15).Value = 9500
$ws.Cells.Item(1184, 16).Value = 9270
$ws.Cells.Item(1184, 19).Value = 464
$ws.Cells.Item(1185, 4).Value = 44707
$ws.Cells.Item(1185, 11).Value = 'Sin especificar'
$ws.Cells.Item(1185, 12).Value = 'Primera Pintón'
$ws.Cells.Item(1185, 13).Value = 3920
$ws.Cells.Item(1185, 14).Value = 9000
$ws.Cells.Item(1185, 15).Value = 10000
$ws.Cells.Item(1185, 16).Value = 9551
$ws.Cells.Item(1185, 19).Value = 478
$ws.Cells.Item(1186, 4).Value = 44265
$ws.Cells.Item(1186, 13).Value = 1120
$ws.Cells.Item(1186, 14).Value = 11000
$ws.Cells.Item(1186, 15).Value = 12000
$ws.Cells.Item(1186, 16).Value = 11500
$ws.Cells.Item(1186, 19).Value = 575
$ws.Cells.Item(1187, 4).Value = 44265
$ws.Cells.Item(1187, 13).Value = 1200
$ws.Cells.Item(1187, 14).Value = 13000
$ws.Cells.Item(1187, 16).Value = 13500
$ws.Cells.Item(1187, 19).Value = 675
$ws.Cells.Item(1188, 4).Value = 44421
$ws.Cells.Item(1188, 11).Value = 'Barraganete'
$ws.Cells.Item(1188, 12).Value = 'Primera'
$ws.Cells.Item(1188, 13).Value = 216
$ws.Cells.Item(1188, 14).Value = 20000
$ws.Cells.Item(1188, 15).Value = 21000
$ws.Cells.Item(1188, 16).Value = 20500
$ws.Cells.Item(1188, 19).Value = 1025
$ws.Cells.Item(1189, 4).Value = 44421
$ws.Cells.Item(1189, 12).Value = 'Pintón'
$ws.Cells.Item(1189, 13).Value = 1280
$ws.Cells.Item(1189, 15).Value = 11000
$ws.Cells.Item(1189, 16).Value = 10688
$ws.Cells.Item(1189, 19).Value = 534
$ws.Cells.Item(1190, 4).Value = 44421
$ws.Cells.Item(1190, 12).Value = 'Primera Pintón'
$ws.Cells.Item(1190, 13).Value = 1840
$ws.Cells.Item(1190, 14).Value = 12000
$ws.Cells.Item(1190, 15).Value = 14000
$ws.Cells.Item(1190, 16).Value = 12891
$ws.Cells.Item(1190, 19).Value = 645
$ws.Cells.Item(1191, 4).Value = 44215
$ws.Cells.Item(1191, 13).Value = 1370
$ws.Cells.Item(1191, 14).Value = 9500
$ws.Cells.Item(1191, 15).Value = 10000
$ws.Cells.Item(1191, 16).Value = 9796
$ws.Cells.Item(1191, 19).Value = 490
$ws.Cells.Item(1192, 4).Value = 44215
$ws.Cells.Item(1192, 12).Value = 'Primera Pintón'
$ws.Cells.Item(1192, 13).Value = 1500
$ws.Cells.Item(1192, 14).Value = 10000
$ws.Cells.Item(1192, 15).Value = 13000
$ws.Cells.Item(1192, 16).Value = 11400
$ws.Cells.Item(1192, 19).Value = 570
$ws.Cells.Item(1193, 12).Value = 'Maduro'
$ws.Cells.Item(1193, 13).Value = 1940
$ws.Cells.Item(1193, 14).Value = 8000
$ws.Cells.Item(1193, 15).Value = 9000
$ws.Cells.Item(1193, 16).Value = 8649
$ws.Cells.Item(1193, 19).Value = 432
$ws.Cells.Item(1194, 4).Value = 44566
$ws.Cells.Item(1194, 12).Value = 'Pintón'
$ws.Cells.Item(1194, 14).Value = 10000
$ws.Cells.Item(1194, 15).Value = 11000
$ws.Cells.Item(1194, 16).Value = 10500
$ws.Cells.Item(1194, 19).Value = 525
$ws.Cells.Item(1195, 4).Value = 44566
$ws.Cells.Item(1195, 12).Value = 'Primera Maduro'
$ws.Cells.Item(1195, 13).Value = 920
$ws.Cells.Item(1195, 14).Value = 11000
$ws.Cells.Item(1195, 15).Value = 12000
$ws.Cells.Item(1195, 16).Value = 11565
$ws.Cells.Item(1195, 19).Value = 578
$ws.Cells.Item(1196, 4).Value = 44566
$ws.Cells.Item(1196, 12).Value = 'Primera Pintón'
$ws.Cells.Item(1196, 13).Value = 3280
$ws.Cells.Item(1196, 15).Value = 13000
$ws.Cells.Item(1196, 16).Value = 11890
$ws.Cells.Item(1196, 19).Value = 594
$ws.Cells.Item(1197, 4).Value = 44637
$ws.Cells.Item(1197, 13).Value = 960
$ws.Cells.Item(1197, 14).Value = 18000
$ws.Cells.Item(1197, 15).Value = 19000
$ws.Cells.Item(1197, 16).Value = 18583
$ws.Cells.Item(1197, 19).Value = 929
$ws.Cells.Item(1198, 4).Value = 44637
$ws.Cells.Item(1198, 13).Value = 1040
$ws.Cells.Item(1198, 14).Value = 19000
$ws.Cells.Item(1198, 15).Value = 20000
$ws.Cells.Item(1198, 16).Value = 19538
$ws.Cells.Item(1198, 19).Value = 977
$ws.Cells.Item(1199, 12).Value = 'Pintón'
$ws.Cells.Item(1199, 13).Value = 2480
$ws.Cells.Item(1199, 14).Value = 11000
$ws.Cells.Item(1199, 15).Value = 14000
$ws.Cells.Item(1199, 16).Value = 12435
$ws.Cells.Item(1199, 19).Value = 622
$ws.Cells.Item(1200, 12).Value = 'Primera Maduro'
$ws.Cells.Item(1200, 13).Value = 840
$ws.Cells.Item(1200, 14).Value = 11000
$ws.Cells.Item(1200, 15).Value = 12000
$ws.Cells.Item(1200, 16).Value = 11571
$ws.Cells.Item(1200, 19).Value = 579
$ws.Cells.Item(1201, 4).Value = 44195
$ws.Cells.Item(1201, 12).Value = 'Primera Pintón'
$ws.Cells.Item(1201, 13).Value = 3280
$ws.Cells.Item(1201, 14).Value = 12500
$ws.Cells.Item(1201, 15).Value = 15000
$ws.Cells.Item(1201, 16).Value = 13591
$ws.Cells.Item(1201, 19).Value = 680
$ws.Cells.Item(1202, 4).Value = 44195
$ws.Cells.Item(1202, 12).Value = 'Primera Verde'
$ws.Cells.Item(1202, 13).Value = 560
$ws.Cells.Item(1202, 14).Value = 13500
$ws.Cells.Item(1202, 15).Value = 13500
$ws.Cells.Item(1202, 16).Value = 13500
$ws.Cells.Item(1202, 19).Value = 675
$ws.Cells.Item(1203, 4).Value = 44195
$ws.Cells.Item(1203, 12).Value = 'Verde'
$ws.Cells.Item(1203, 13).Value = 400
$ws.Cells.Item(1203, 14).Value = 13000
$ws.Cells.Item(1203, 16).Value = 13000
$ws.Cells.Item(1203, 19).Value = 650
$ws.Cells.Item(1204, 4).Value = 44244
$ws.Cells.Item(1204, 11).Value = 'Sin especificar'
$ws.Cells.Item(1204, 12).Value = 'Pintón'
$ws.Cells.Item(1204, 13).Value = 1490
$ws.Cells.Item(1204, 14).Value = 9000
$ws.Cells.Item(1204, 15).Value = 10000
$ws.Cells.Item(1204, 16).Value = 9638
$ws.Cells.Item(1204, 19).Value = 482
$ws.Cells.Item(1205, 4).Value = 44244
$ws.Cells.Item(1205, 12).Value = 'Primera Maduro'
$ws.Cells.Item(1205, 13).Value = 760
$ws.Cells.Item(1205, 14).Value = 11000
$ws.Cells.Item(1205, 16).Value = 11526
$ws.Cells.Item(1205, 19).Value = 576
$ws.Cells.Item(1206, 4).Value = 44244
$ws.Cells.Item(1206, 12).Value = 'Primera Pintón'
$ws.Cells.Item(1206, 13).Value = 2470
$ws.Cells.Item(1206, 14).Value = 12000
$ws.Cells.Item(1206, 15).Value = 13000
$ws.Cells.Item(1206, 16).Value = 12324
$ws.Cells.Item(1206, 19).Value = 616
$ws.Cells.Item(1207, 11).Value = 'Barraganete'
$ws.Cells.Item(1207, 12).Value = 'Primera'
$ws.Cells.Item(1207, 13).Value = 316
$ws.Cells.Item(1207, 14).Value = 19000
$ws.Cells.Item(1207, 15).Value = 20000
$ws.Cells.Item(1207, 16).Value = 19658
$ws.Cells.Item(1207, 19).Value = 983
$ws.Cells.Item(1208, 4).Value = 44442
$ws.Cells.Item(1208, 13).Value = 150
$ws.Cells.Item(1208, 14).Value = 12000
$ws.Cells.Item(1208, 15).Value = 12000
$ws.Cells.Item(1208, 16).Value = 12000
$ws.Cells.Item(1208, 19).Value = 600
$ws.Cells.Item(1209, 4).Value = 44442
$ws.Cells.Item(1209, 13).Value = 1360
$ws.Cells.Item(1209, 14).Value = 14000
$ws.Cells.Item(1209, 15).Value = 16000
$ws.Cells.Item(1209, 16).Value = 15118
$ws.Cells.Item(1209, 19).Value = 756
$ws.Cells.Item(1210, 4).Value = 44442
$ws.Cells.Item(1210, 13).Value = 2015
$ws.Cells.Item(1210, 14).Value = 16000
$ws.Cells.Item(1210, 15).Value = 18000
$ws.Cells.Item(1210, 16).Value = 17079
$ws.Cells.Item(1210, 19).Value = 854
$ws.Cells.Item(1211, 12).Value = 'Maduro'
$ws.Cells.Item(1211, 13).Value = 820
$ws.Cells.Item(1211, 14).Value = 17000
$ws.Cells.Item(1211, 15).Value = 18000
$ws.Cells.Item(1211, 16).Value = 17463
$ws.Cells.Item(1211, 19).Value = 873
$ws.Cells.Item(1212, 4).Value = 44483
$ws.Cells.Item(1212, 12).Value = 'Pintón'
$ws.Cells.Item(1212, 13).Value = 575
$ws.Cells.Item(1212, 14).Value = 18000
$ws.Cells.Item(1212, 15).Value = 22000
$ws.Cells.Item(1212, 16).Value = 19217
$ws.Cells.Item(1212, 19).Value = 961
$ws.Cells.Item(1213, 4).Value = 44483
$ws.Cells.Item(1213, 12).Value = 'Primera Pintón'
$ws.Cells.Item(1213, 13).Value = 2455
$ws.Cells.Item(1213, 14).Value = 19000
$ws.Cells.Item(1213, 15).Value = 23000
$ws.Cells.Item(1213, 16).Value = 20703
$ws.Cells.Item(1213, 19).Value = 1035
$ws.Cells.Item(1214, 4).Value = 44483
$ws.Cells.Item(1214, 12).Value = 'Primera Verde'
$ws.Cells.Item(1214, 13).Value = 1820
$ws.Cells.Item(1214, 14).Value = 21000
$ws.Cells.Item(1214, 15).Value = 22000
$ws.Cells.Item(1214, 16).Value = 21516
$ws.Cells.Item(1214, 19).Value = 1076
$ws.Cells.Item(1215, 12).Value = 'Maduro'
$ws.Cells.Item(1215, 13).Value = 995
$ws.Cells.Item(1215, 14).Value = 13000
$ws.Cells.Item(1215, 15).Value = 15000
$ws.Cells.Item(1215, 16).Value = 14246
$ws.Cells.Item(1215, 19).Value = 712
$ws.Cells.Item(1216, 4).Value = 44663
$ws.Cells.Item(1216, 13).Value = 1295
$ws.Cells.Item(1216, 14).Value = 14000
$ws.Cells.Item(1216, 15).Value = 15000
$ws.Cells.Item(1216, 16).Value = 14598
$ws.Cells.Item(1216, 19).Value = 730
$ws.Cells.Item(1217, 4).Value = 44663
$ws.Cells.Item(1217, 12).Value = 'Primera Maduro'
$ws.Cells.Item(1217, 13).Value = 1040
$ws.Cells.Item(1217, 14).Value = 16000
$ws.Cells.Item(1217, 15).Value = 17000
$ws.Cells.Item(1217, 16).Value = 16558
$ws.Cells.Item(1217, 19).Value = 828
$ws.Cells.Item(1218, 4).Value = 44663
$ws.Cells.Item(1218, 12).Value = 'Primera Pintón'
$ws.Cells.Item(1218, 13).Value = 3575
$ws.Cells.Item(1218, 14).Value = 16000
$ws.Cells.Item(1218, 15).Value = 17000
$ws.Cells.Item(1218, 16).Value = 16392
$ws.Cells.Item(1218, 19).Value = 820
$ws.Cells.Item(1219, 4).Value = 44188
$ws.Cells.Item(1219, 12).Value = 'Pintón'
$ws.Cells.Item(1219, 13).Value = 2010
$ws.Cells.Item(1219, 14).Value = 9000
$ws.Cells.Item(1219, 15).Value = 10000
$ws.Cells.Item(1219, 16).Value = 9562
$ws.Cells.Item(1219, 19).Value = 478
$ws.Cells.Item(1220, 4).Value = 44188
$ws.Cells.Item(1220, 12).Value = 'Primera Pintón'
$ws.Cells.Item(1220, 13).Value = 2020
$ws.Cells.Item(1220, 14).Value = 11000
$ws.Cells.Item(1220, 15).Value = 12000
$ws.Cells.Item(1220, 16).Value = 11574
$ws.Cells.Item(1220, 19).Value = 579
$ws.Cells.Item(1221, 4).Value = 44187
$ws.Cells.Item(1221, 12).Value = 'Pintón'
$ws.Cells.Item(1221, 13).Value = 1520
$ws.Cells.Item(1221, 14).Value = 9000
$ws.Cells.Item(1221, 15).Value = 10000
$ws.Cells.Item(1221, 16).Value = 9500
$ws.Cells.Item(1221, 19).Value = 475
$ws.Cells.Item(1222, 4).Value = 44187
$ws.Cells.Item(1222, 13).Value = 3200
$ws.Cells.Item(1222, 14).Value = 10000
$ws.Cells.Item(1222, 15).Value = 12000
$ws.Cells.Item(1222, 16).Value = 10875
$ws.Cells.Item(1222, 19).Value = 544
$ws.Cells.Item(1223, 4).Value = 44519
$ws.Cells.Item(1223, 12).Value = 'Pintón'
$ws.Cells.Item(1223, 13).Value = 1160
$ws.Cells.Item(1223, 14).Value = 13000
$ws.Cells.Item(1223, 15).Value = 14000
$ws.Cells.Item(1223, 16).Value = 13414
$ws.Cells.Item(1223, 19).Value = 671
$ws.Cells.Item(1224, 4).Value = 44519
$ws.Cells.Item(1224, 12).Value = 'Primera Maduro'
$ws.Cells.Item(1224, 13).Value = 640
$ws.Cells.Item(1224, 14).Value = 13000
$ws.Cells.Item(1224, 15).Value = 14000
$ws.Cells.Item(1224, 16).Value = 13562
$ws.Cells.Item(1224, 19).Value = 678
$ws.Cells.Item(1225, 4).Value = 44519
$ws.Cells.Item(1225, 13).Value = 3540
$ws.Cells.Item(1225, 15).Value = 15000
$ws.Cells.Item(1225, 16).Value = 14486
$ws.Cells.Item(1225, 19).Value = 724
$ws.Cells.Item(1226, 4).Value = 44231
$ws.Cells.Item(1226, 12).Value = 'Maduro'
$ws.Cells.Item(1226, 13).Value = 200
$ws.Cells.Item(1226, 14).Value = 10000
$ws.Cells.Item(1226, 15).Value = 10000
$ws.Cells.Item(1226, 16).Value = 10000
$ws.Cells.Item(1226, 19).Value = 500
$ws.Cells.Item(1227, 4).Value = 44231
$ws.Cells.Item(1227, 12).Value = 'Pintón'
$ws.Cells.Item(1227, 13).Value = 200
$ws.Cells.Item(1227, 15).Value = 12000
$ws.Cells.Item(1227, 16).Value = 12000
$ws.Cells.Item(1227, 19).Value = 600
$ws.Cells.Item(1228, 4).Value = 44231
$ws.Cells.Item(1228, 13).Value = 2740
$ws.Cells.Item(1228, 14).Value = 14000
$ws.Cells.Item(1228, 15).Value = 16000
$ws.Cells.Item(1228, 16).Value = 15234
$ws.Cells.Item(1228, 19).Value = 762
$ws.Cells.Item(1229, 1).Value = 6
$ws.Cells.Item(1229, 2).Value = 'Mercado Mayorista Lo Valledor de Santiago'
$ws.Cells.Item(1229, 3).Value = 'Metropolitana'
$ws.Cells.Item(1229, 4).Value = 44194
$ws.Cells.Item(1229, 5).Value = 13
$ws.Cells.Item(1229, 6).Value = 'Fruta'
$ws.Cells.Item(1229, 7).Value = 100108
$ws.Cells.Item(1229, 8).Value = 'Tropicales y subtropicales'
$ws.Cells.Item(1229, 9).Value = 100108006
$ws.Cells.Item(1229, 10).Value = 'Plátano'
$ws.Cells.Item(1229, 11).Value = 'Sin especificar'
$ws.Cells.Item(1229, 12).Value = 'Pintón'
$ws.Cells.Item(1229, 13).Value = 2080
$ws.Cells.Item(1229, 14).Value = 11000
$ws.Cells.Item(1229, 15).Value = 12000
$ws.Cells.Item(1229, 16).Value = 11500
$ws.Cells.Item(1229, 17).Value = '$/caja 20 kilos'
$ws.Cells.Item(1229, 18).Value = 'Ecuador'
$ws.Cells.Item(1229, 19).Value = 575
$ws.Cells.Item(1229, 20).Value = 20
$ws.Cells.Item(1230, 1).Value = 6
$ws.Cells.Item(1230, 2).Value = 'Mercado Mayorista Lo Valledor de Santiago'
$ws.Cells.Item(1230, 3).Value = 'Metropolitana'
$ws.Cells.Item(1230, 4).Value = 44194
$ws.Cells.Item(1230, 5).Value = 13
$ws.Cells.Item(1230, 6).Value = 'Fruta'
$ws.Cells.Item(1230, 7).Value = 100108
$ws.Cells.Item(1230, 8).Value = 'Tropicales y subtropicales'
$ws.Cells.Item(1230, 9).Value = 100108006
$ws.Cells.Item(1230, 10).Value = 'Plátano'
$ws.Cells.Item(1230, 11).Value = 'Sin especificar'
$ws.Cells.Item(1230, 12).Value = 'Primera Maduro'
$ws.Cells.Item(1230, 13).Value = 960
$ws.Cells.Item(1230, 14).Value = 12000
$ws.Cells.Item(1230, 15).Value = 13000
$ws.Cells.Item(1230, 16).Value = 12417
$ws.Cells.Item(1230, 17).Value = '$/caja 20 kilos'
$ws.Cells.Item(1230, 18).Value = 'Ecuador'
$ws.Cells.Item(1230, 19).Value = 621
$ws.Cells.Item(1230, 20).Value = 20
$ws.Cells.Item(1231, 1).Value = 6
$ws.Cells.Item(1231, 2).Value = 'Mercado Mayorista Lo Valledor de Santiago'
$ws.Cells.Item(1231, 3).Value = 'Metropolitana'
$ws.Cells.Item(1231, 4).Value = 44194
$ws.Cells.Item(1231, 5).Value = 13
$ws.Cells.Item(1231, 6).Value = 'Fruta'
$ws.Cells.Item(1231, 7).Value = 100108
$ws.Cells.Item(1231, 8).Value = 'Tropicales y subtropicales'
$ws.Cells.Item(1231, 9).Value = 100108006
$ws.Cells.Item(1231, 10).Value = 'Plátano'
$ws.Cells.Item(1231, 11).Value = 'Sin especificar'
$ws.Cells.Item(1231, 12).Value = 'Primera Pintón'
$ws.Cells.Item(1231, 13).Value = 3340
$ws.Cells.Item(1231, 14).Value = 12000
$ws.Cells.Item(1231, 15).Value = 14000
$ws.Cells.Item(1231, 16).Value = 13165
$ws.Cells.Item(1231, 17).Value = '$/caja 20 kilos'
$ws.Cells.Item(1231, 18).Value = 'Ecuador'
$ws.Cells.Item(1231, 19).Value = 658
$ws.Cells.Item(1231, 20).Value = 20

# Preserve the date number format on the newly created rows' Fecha column (D)
$dateFormat = $ws.Range("D1228").NumberFormat
$ws.Range("D1229").NumberFormat = $dateFormat
$ws.Range("D1230").NumberFormat = $dateFormat
$ws.Range("D1231").NumberFormat = $dateFormat
